# Convert one column of an xlsx file to txt file
# - The value in C3 ("Chennai") is updated to "Bangalore".
# - The active/selected cell moves to C3 (reflecting the edit location).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Bangalore"
$ws.Range("C3").Select()
